$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.276.52"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.878.97"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.19%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.13"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.69%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2879"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06591"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.904.65"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.87"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.14%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.183"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.81"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6606"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.255.89"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.46"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007725"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.448"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.148.36"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "192.25"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.190"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.432"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.95"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.41"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.941"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.444"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.262"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.60%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.050"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05078"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7449"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.97%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01834"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.634"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9150"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.074"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.47"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.886"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4328"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9996"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.683"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1357"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.579"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +9.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.23"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -9.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.951"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05730"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.65%  "
